# Add customer 345678 to mock data (new row 4 on Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-F and H hold text-like values ("345678", ISO dates, words).
# Pre-format as Text so Excel's auto-detection doesn't turn the
# numeric-looking id or the date-looking strings into a number / date
# serial, matching the literal strings stored in the other data rows.
$ws.Range("A4:F4").NumberFormat = "@"
$ws.Range("H4").NumberFormat = "@"

$ws.Range("A4").Value = "345678"
$ws.Range("B4").Value = "2012-08-10"
$ws.Range("C4").Value = "2012-09-01"
$ws.Range("D4").Value = "Employed"
$ws.Range("E4").Value = "Male"
$ws.Range("F4").Value = "Married"
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = "India"

# Restore the default (unstyled) look so the new row matches rows 2-3,
# which carry no explicit cell style.
$ws.Range("A4:F4").Style = "Normal"
$ws.Range("H4").Style = "Normal"
